$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Feria Lagunitas de Puerto Montt -
# Poroto granado". It slots in at row 20, pushing the previously existing
# rows 20-62 down to rows 21-63 (Excel's normal "insert row" shift).
$ws.Rows.Item(20).Insert()

# Fill in the data for the newly inserted row 20.
$ws.Cells.Item(20, 1).Value = 4
$ws.Cells.Item(20, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(20, 3).Value = "Los Lagos"
$ws.Cells.Item(20, 4).Value = 44967
$ws.Cells.Item(20, 5).Value = 10
$ws.Cells.Item(20, 6).Value = 100112030
$ws.Cells.Item(20, 7).Value = "Poroto granado"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 80
$ws.Cells.Item(20, 11).Value = 40000
$ws.Cells.Item(20, 12).Value = 40000
$ws.Cells.Item(20, 13).Value = 40000
$ws.Cells.Item(20, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(20, 15).Value = "Región Metropolitana"
$ws.Cells.Item(20, 16).Value = 1600
$ws.Cells.Item(20, 17).Value = 25
$ws.Cells.Item(20, 18).Value = "Hortaliza"
